$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Authors column (E) for the two reference rows with refreshed
# author-list strings (elastic search json query output produced an extra
# leading-space increment versus the previous extraction pass).
$ws.Range("E2").Value = "[Anthony V%Das%NULL%1,                           Padmaja K%Rani%NULL%1,                           Pravin K%Vaddavalli%NULL%1]"
$ws.Range("E3").Value = "[Gagan%Kalra%NULL%1,                           Andrew M.%Williams%NULL%1,                           Patrick W.%Commiskey%NULL%1,                           Eve M. R.%Bowers%NULL%1,                           Tadhg%Schempf%NULL%1,                           José-Alain%Sahel%NULL%1,                           Evan L.%Waxman%waxmane@upmc.edu%1,                           Roxana%Fu%fur3@upmc.edu%1]"
